$p = $ppt.ActivePresentation

$oldDate = "1/12/2022"
$newDate = "9/27/2023"

function Update-DatePlaceholder {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 1 subtitle: "January 12, 2022" -> "University of Mount Union".
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shape = $slide1.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "January 12, 2022") {
            $tr.Text = "University of Mount Union"
        }
    }
}
